$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Q7)
$ws.Range("B9").Value = 0.0315757522548537
$ws.Range("C9").Value = 0.3042284657534167
$ws.Range("D9").Value = 0.1264694874506458
$ws.Range("E9").Value = 0.3556254876279901
$ws.Range("F9").Value = 0.388029575141152
$ws.Range("G9").Value = 6

# Row 10 (Q8)
$ws.Range("B10").Value = 0.1190308503644765
$ws.Range("C10").Value = 0.1190308503644765
$ws.Range("D10").Value = 0.02079770754726164
$ws.Range("E10").Value = 0.144214103149663
$ws.Range("F10").Value = 0.09971983911517755
$ws.Range("G10").Value = 3

# Row 11 (Q9)
$ws.Range("B11").Value = -0.04243697084963852
$ws.Range("C11").Value = 0.04243697084963852
$ws.Range("D11").Value = 0.001800896494893069
$ws.Range("E11").Value = 0.04243697084963852
$ws.Range("F11").ClearContents()
$ws.Range("G11").Value = 1
